$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row whose values/formatting are replicated into the new rows.
$srcRow = 89
$numNewRows = 3

for ($i = 1; $i -le $numNewRows; $i++) {
    $destRow = $srcRow + $i

    # Copy formatting only (keeps the date style s="2" on column A, plain
    # numeric style on the rest) from the last existing row into the new row.
    $ws.Range("A" + $srcRow + ":J" + $srcRow).Copy()
    $ws.Range("A" + $destRow + ":J" + $destRow).PasteSpecial(-4122)

    # Column A holds consecutive daily dates, so bump it by $i days relative
    # to the source row. The rest of the columns repeat the source row's values.
    $ws.Cells.Item($destRow, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2 + $i
    for ($c = 2; $c -le 10; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $ws.Cells.Item($srcRow, $c).Value2
    }
}

$excel.CutCopyMode = $false
